$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 34.582962
$ws.Range("H2").Value = 103.748886
$ws.Range("I2").Value = 0.8305266248561333
$ws.Range("J2").Value = 0.8305266248561333
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.324764666666667
$ws.Range("N2").Value = 6.974294
$ws.Range("O2").Value = 0.04473923998638302
$ws.Range("P2").Value = 0.04473923998638301
$ws.Range("Q2").Value = 80.39724812627601
$ws.Range("R2").Value = 723.575233136484
$ws.Range("S2").Value = 0.03715712998451925
$ws.Range("T2").Value = 0.03715712998451924

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 34.582962
$ws.Range("H3").Value = 103.748886
$ws.Range("I3").Value = 0.8305266248561333
$ws.Range("J3").Value = 0.8305266248561333
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 19.27491966666667
$ws.Range("N3").Value = 57.824759
$ws.Range("O3").Value = 0.3709387315842666
$ws.Range("P3").Value = 0.3709387315842665
$ws.Range("Q3").Value = 666.583814385386
$ws.Range("R3").Value = 5999.254329468474
$ws.Range("S3").Value = 0.3080744927710961
$ws.Range("T3").Value = 0.3080744927710961

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 34.582962
$ws.Range("H4").Value = 103.748886
$ws.Range("I4").Value = 0.8305266248561333
$ws.Range("J4").Value = 0.8305266248561333
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 30.36285833333334
$ws.Range("N4").Value = 91.088575
$ws.Range("O4").Value = 0.5843220284293504
$ws.Range("P4").Value = 0.5843220284293504
$ws.Range("Q4").Value = 1050.03757595305
$ws.Range("R4").Value = 9450.33818357745
$ws.Range("S4").Value = 0.485295002100518
$ws.Range("T4").Value = 0.485295002100518

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.670847666666667
$ws.Range("H5").Value = 14.012543
$ws.Range("I5").Value = 0.1121726747354322
$ws.Range("J5").Value = 0.1121726747354322
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.324764666666667
$ws.Range("N5").Value = 6.974294
$ws.Range("O5").Value = 0.04473923998638302
$ws.Range("P5").Value = 0.04473923998638301
$ws.Range("Q5").Value = 10.85862161884911
$ws.Range("R5").Value = 97.72759456964201
$ws.Range("S5").Value = 0.005018520214902985
$ws.Range("T5").Value = 0.005018520214902985

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 4.670847666666667
$ws.Range("H6").Value = 14.012543
$ws.Range("I6").Value = 0.1121726747354322
$ws.Range("J6").Value = 0.1121726747354322
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 19.27491966666667
$ws.Range("N6").Value = 57.824759
$ws.Range("O6").Value = 0.3709387315842666
$ws.Range("P6").Value = 0.3709387315842665
$ws.Range("Q6").Value = 90.03021355023745
$ws.Range("R6").Value = 810.2719219521371
$ws.Range("S6").Value = 0.04160918968477573
$ws.Range("T6").Value = 0.04160918968477573

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4.670847666666667
$ws.Range("H7").Value = 14.012543
$ws.Range("I7").Value = 0.1121726747354322
$ws.Range("J7").Value = 0.1121726747354322
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 30.36285833333334
$ws.Range("N7").Value = 91.088575
$ws.Range("O7").Value = 0.5843220284293504
$ws.Range("P7").Value = 0.5843220284293504
$ws.Range("Q7").Value = 141.8202859995806
$ws.Range("R7").Value = 1276.382573996225
$ws.Range("S7").Value = 0.06554496483575349
$ws.Range("T7").Value = 0.0655449648357535

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.385989666666667
$ws.Range("H8").Value = 7.157969
$ws.Range("I8").Value = 0.05730070040843456
$ws.Range("J8").Value = 0.05730070040843457
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.324764666666667
$ws.Range("N8").Value = 6.974294
$ws.Range("O8").Value = 0.04473923998638302
$ws.Range("P8").Value = 0.04473923998638301
$ws.Range("Q8").Value = 5.546864472098445
$ws.Range("R8").Value = 49.921780248886
$ws.Range("S8").Value = 0.002563589786960789
$ws.Range("T8").Value = 0.002563589786960789

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.385989666666667
$ws.Range("H9").Value = 7.157969
$ws.Range("I9").Value = 0.05730070040843456
$ws.Range("J9").Value = 0.05730070040843457
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 19.27491966666667
$ws.Range("N9").Value = 57.824759
$ws.Range("O9").Value = 0.3709387315842666
$ws.Range("P9").Value = 0.3709387315842665
$ws.Range("Q9").Value = 45.98975915049677
$ws.Range("R9").Value = 413.907832354471
$ws.Range("S9").Value = 0.02125504912839478
$ws.Range("T9").Value = 0.02125504912839478

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.385989666666667
$ws.Range("H10").Value = 7.157969
$ws.Range("I10").Value = 0.05730070040843456
$ws.Range("J10").Value = 0.05730070040843457
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 30.36285833333334
$ws.Range("N10").Value = 91.088575
$ws.Range("O10").Value = 0.5843220284293504
$ws.Range("P10").Value = 0.5843220284293504
$ws.Range("Q10").Value = 72.44546623379723
$ws.Range("R10").Value = 652.009196104175
$ws.Range("S10").Value = 0.03348206149307899
$ws.Range("T10").Value = 0.03348206149307899
